$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.161.57'
$ws.Range("E2").Value = '  -3.89%  '
$ws.Range("D3").Value = '2.456.53'
$ws.Range("E3").Value = '  -3.30%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'310.55"
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("D6").Value = "'94.52"
$ws.Range("E6").Value = '  -6.93%  '
$ws.Range("D7").Value = "'0.550"
$ws.Range("E7").Value = '  -3.87%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = "'0.505"
$ws.Range("E9").Value = '  -4.61%  '
$ws.Range("D10").Value = "'33.52"
$ws.Range("E10").Value = '  -7.79%  '
$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = '  -2.94%  '
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("D13").Value = "'6.97"
$ws.Range("E13").Value = '  -5.32%  '
$ws.Range("D14").Value = '2.829.41'
$ws.Range("E14").Value = '  -3.34%  '
$ws.Range("D15").Value = '2.441.34'
$ws.Range("E15").Value = '  -2.50%  '
$ws.Range("D16").Value = "'14.46"
$ws.Range("E16").Value = '  -8.68%  '
$ws.Range("D17").Value = "'0.787"
$ws.Range("E17").Value = '  -3.45%  '
$ws.Range("D18").Value = '41.125.75'
$ws.Range("E18").Value = '  -3.91%  '
$ws.Range("D19").Value = "'6.36"
$ws.Range("E19").Value = '  -6.30%  '
$ws.Range("D20").Value = '0.0₃0915'
$ws.Range("E20").Value = '  -4.10%  '
$ws.Range("D21").Value = "'11.56"
$ws.Range("E21").Value = '  -6.09%  '
$ws.Range("D22").Value = "'67.24"
$ws.Range("E22").Value = '  -3.06%  '
$ws.Range("D23").Value = "'237.76"
$ws.Range("E23").Value = '  -2.81%  '
$ws.Range("D24").Value = "'2.77"
$ws.Range("E24").Value = '  -4.85%  '
$ws.Range("D25").Value = "'1.94"
$ws.Range("E25").Value = '  -5.73%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").Value = "'24.61"
$ws.Range("E27").Value = '  -5.93%  '
$ws.Range("E28").Value = '  -3.85%  '
$ws.Range("D29").Value = "'9.68"
$ws.Range("E29").Value = '  -5.31%  '
$ws.Range("D30").Value = "'36.24"
$ws.Range("E30").Value = '  -8.48%  '
$ws.Range("D31").Value = "'152.94"
$ws.Range("D32").Value = "'5.60"
$ws.Range("E32").Value = '  -3.62%  '
$ws.Range("E33").Value = '  -1.11%  '
$ws.Range("D34").Value = "'0.0752"
$ws.Range("E34").Value = '  -5.46%  '
$ws.Range("E35").Value = '  -8.76%  '
$ws.Range("D36").Value = "'3.03"
$ws.Range("E36").Value = '  -5.16%  '
$ws.Range("D37").Value = "'17.28"
$ws.Range("E37").Value = '  -6.21%  '
$ws.Range("E38").Value = '  -7.26%  '
$ws.Range("D39").Value = "'0.104"
$ws.Range("E39").Value = '  -7.75%  '
$ws.Range("D40").Value = "'0.114"
$ws.Range("E40").Value = '  -4.34%  '
$ws.Range("D41").Value = "'4.19"
$ws.Range("E41").Value = '  -3.07%  '
$ws.Range("D42").Value = "'21.34"
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '1.961.86'
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("D45").Value = "'0.0284"
$ws.Range("E45").Value = '  -5.10%  '
$ws.Range("D46").Value = "'3.05"
$ws.Range("E46").Value = '  -8.56%  '
$ws.Range("D47").Value = "'8.72"
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("D48").Value = "'77.00"
$ws.Range("E48").Value = '  -4.84%  '
$ws.Range("D49").Value = "'97.56"
$ws.Range("E49").Value = '  -3.52%  '
$ws.Range("D50").Value = "'69.14"
$ws.Range("E50").Value = '  -5.07%  '
$ws.Range("D51").Value = "'0.180"
$ws.Range("E51").Value = '  -6.91%  '
